$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.669
$ws.Range("E3").Value = 12.457

$ws.Range("C6").Value = -12.974

$ws.Range("C7").Value = -13.057

$ws.Range("B8").Value = 5.991999999999999
$ws.Range("C8").Value = -12.194
$ws.Range("E8").Value = 12.711

$ws.Range("D11").Value = -8.334

$ws.Range("A12").Value = -21.482
$ws.Range("B12").Value = 6.806999999999999

$ws.Range("B14").Value = 6.973000000000001
$ws.Range("D14").Value = -7.764

$ws.Range("C19").Value = -12.912
$ws.Range("D19").Value = -7.928
$ws.Range("E19").Value = 12.806

$ws.Range("C21").Value = -12.657
$ws.Range("D21").Value = -7.528

$ws.Range("B22").Value = 6.74

$ws.Range("C24").Value = -12.512
